$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated daily hourly spot prices for row 2 (automatic electricity price update)
$ws.Range("A2").Value = 45919
$ws.Range("B2").Value = 113.07
$ws.Range("C2").Value = 110.91
$ws.Range("D2").Value = 108.34
$ws.Range("E2").Value = 108.26
$ws.Range("F2").Value = 107.21
$ws.Range("G2").Value = 108.26
$ws.Range("H2").Value = 110.78
$ws.Range("I2").Value = 120.42
$ws.Range("J2").Value = 113.54
$ws.Range("K2").Value = 106.68
$ws.Range("L2").Value = 66.04000000000001
$ws.Range("M2").Value = 52.14
$ws.Range("N2").Value = 49
$ws.Range("O2").Value = 45.1
$ws.Range("P2").Value = 31.25
$ws.Range("Q2").Value = 32.61
$ws.Range("R2").Value = 55
$ws.Range("S2").Value = 65.87
$ws.Range("T2").Value = 76
$ws.Range("U2").Value = 108.14
$ws.Range("V2").Value = 125.2
$ws.Range("W2").Value = 125
$ws.Range("X2").Value = 105.78
$ws.Range("Y2").Value = 104.68
$ws.Range("Z2").Value = 89.55
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 115.16
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 125.1
$ws.Range("AE2").Value = "6h-8h"
$ws.Range("AF2").Value = 115.6
$ws.Range("AG2").Value = "10h-18h"
